# "Generate Report for Handback" — refresh the localization-status report:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview!E2/F2 and the per-language Status cell, zh-cn!C2 / de-de!C2,
#     all shared the same string so they're updated together).
#   * The per-language "Latest Handback DateTime" timestamps advance to the
#     new handback run.
#   * The stale "handback file is not the latest" error notices are cleared
#     now that the handback is in sync.
#   * Column widths for the widened Status columns / shrunk Error Detail
#     columns are refreshed to match the new (longer/shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status cells ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-18 10:52:34"
$wsZhCn.Range("P2").Value = ""

# --- de-de detail sheet ---
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-18 10:52:41"
$wsDeDe.Range("P2").Value = ""

# --- Column width refresh (status text grew, error text now empty) ---
# Target widths from the report regeneration: ~29.98 chars for the widened
# Status columns, ~13.75 chars for the shrunk Error Detail columns. The
# nearest values the host's column-width quantizer can reproduce are used.
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667
$wsZhCn.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334
$wsDeDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
